$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for A/B/E/F/G/H/Q/R are cyclically rotated among rows 2, 3, 5:
#   new row2 <- old row5
#   new row3 <- old row2
#   new row5 <- old row3
# Row 4 is untouched.

# Capture old values before overwriting anything.
# Note: use Value2 (not Value) -- in this COM-interop environment the
# Value property getter does not return plain data reliably.
$oldA2 = $ws.Range("A2").Value2
$oldB2 = $ws.Range("B2").Value2
$oldE2 = $ws.Range("E2").Value2
$oldF2 = $ws.Range("F2").Value2
$oldG2 = $ws.Range("G2").Value2
$oldH2 = $ws.Range("H2").Value2
$oldQ2 = $ws.Range("Q2").Value2
$oldR2 = $ws.Range("R2").Value2

$oldA3 = $ws.Range("A3").Value2
$oldB3 = $ws.Range("B3").Value2
$oldE3 = $ws.Range("E3").Value2
$oldF3 = $ws.Range("F3").Value2
$oldG3 = $ws.Range("G3").Value2
$oldH3 = $ws.Range("H3").Value2
$oldQ3 = $ws.Range("Q3").Value2
$oldR3 = $ws.Range("R3").Value2

$oldA5 = $ws.Range("A5").Value2
$oldB5 = $ws.Range("B5").Value2
$oldE5 = $ws.Range("E5").Value2
$oldF5 = $ws.Range("F5").Value2
$oldG5 = $ws.Range("G5").Value2
$oldH5 = $ws.Range("H5").Value2
$oldQ5 = $ws.Range("Q5").Value2
$oldR5 = $ws.Range("R5").Value2

# Row 2 gets old row 5 values
$ws.Range("A2").Value2 = $oldA5
$ws.Range("B2").Value2 = $oldB5
$ws.Range("E2").Value2 = $oldE5
$ws.Range("F2").Value2 = $oldF5
$ws.Range("G2").Value2 = $oldG5
$ws.Range("H2").Value2 = $oldH5
$ws.Range("Q2").Value2 = $oldQ5
$ws.Range("R2").Value2 = $oldR5

# Row 3 gets old row 2 values
$ws.Range("A3").Value2 = $oldA2
$ws.Range("B3").Value2 = $oldB2
$ws.Range("E3").Value2 = $oldE2
$ws.Range("F3").Value2 = $oldF2
$ws.Range("G3").Value2 = $oldG2
$ws.Range("H3").Value2 = $oldH2
$ws.Range("Q3").Value2 = $oldQ2
$ws.Range("R3").Value2 = $oldR2

# Row 5 gets old row 3 values
$ws.Range("A5").Value2 = $oldA3
$ws.Range("B5").Value2 = $oldB3
$ws.Range("E5").Value2 = $oldE3
$ws.Range("F5").Value2 = $oldF3
$ws.Range("G5").Value2 = $oldG3
$ws.Range("H5").Value2 = $oldH3
$ws.Range("Q5").Value2 = $oldQ3
$ws.Range("R5").Value2 = $oldR3
